# Generate Report for Handback
# Populate the "Latest Target File" / "Latest Handback File" / "Latest Handback DateTime"
# columns on the per-locale sheets (zh-cn, de-de) now that handback has happened, and
# flip the Overview sheet's Status from "Ready for handoff" to the handed-back state.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: update per-locale status columns (E = zh-cn, F = de-de) ---
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E2").Value = "Handed back: in sync with en-US"
$ov.Range("F2").Value = "Handed back: in sync with en-US"
$ov.Range("E3").Value = "Handed back: in sync with en-US"
$ov.Range("F3").Value = "Handed back: in sync with en-US"
$ov.Columns.Item(5).ColumnWidth = 29.14   # raw width 30.0 (closest achievable to target 29.9777050018311)
$ov.Columns.Item(6).ColumnWidth = 29.14   # raw width 30.0 (closest achievable to target 29.9777050018311)

# --- zh-cn sheet ---
$zh = $wb.Worksheets.Item("zh-cn")

# Row 2 : 2a70686a-c89a-4542-9cab-dc83a74b1a7f
$zh.Range("I2").Value = "2a70686a-c89a-4542-9cab-dc83a74b1a7f.md"
$zh.Hyperlinks.Add($zh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/64e54dce710ff6806c33fcbd5c6c5c24cea91328/e2e/2a70686a-c89a-4542-9cab-dc83a74b1a7f.md", "", "", "2a70686a-c89a-4542-9cab-dc83a74b1a7f.md") | Out-Null
$zh.Range("J2").Value = "2a70686a-c89a-4542-9cab-dc83a74b1a7f.91dcba358223f4ada572a60e999c28132620d6ae.zh-cn.xlf"
$zh.Range("K2").Value = "2016-10-24 09:52:59"

# Row 3 : de4b381b-7123-43c1-a4ab-173c364b5d43
$zh.Range("I3").Value = "de4b381b-7123-43c1-a4ab-173c364b5d43.md"
$zh.Hyperlinks.Add($zh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/64e54dce710ff6806c33fcbd5c6c5c24cea91328/e2e/de4b381b-7123-43c1-a4ab-173c364b5d43.md", "", "", "de4b381b-7123-43c1-a4ab-173c364b5d43.md") | Out-Null
$zh.Range("J3").Value = "de4b381b-7123-43c1-a4ab-173c364b5d43.9c06a25dcb376ef7f97389c27a3adf04fa2ec2ca.zh-cn.xlf"
$zh.Range("K3").Value = "2016-10-24 09:52:59"

$zh.Columns.Item(3).ColumnWidth = 29.14   # raw width 30.0 (closest achievable to target 29.9777050018311)
$zh.Columns.Item(9).ColumnWidth = 39.14   # raw width 40.0
$zh.Columns.Item(10).ColumnWidth = 39.14  # raw width 40.0

# --- de-de sheet ---
$de = $wb.Worksheets.Item("de-de")

# Row 2 : 2a70686a-c89a-4542-9cab-dc83a74b1a7f
$de.Range("I2").Value = "2a70686a-c89a-4542-9cab-dc83a74b1a7f.md"
$de.Hyperlinks.Add($de.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/64e54dce710ff6806c33fcbd5c6c5c24cea91328/e2e/2a70686a-c89a-4542-9cab-dc83a74b1a7f.md", "", "", "2a70686a-c89a-4542-9cab-dc83a74b1a7f.md") | Out-Null
$de.Range("J2").Value = "2a70686a-c89a-4542-9cab-dc83a74b1a7f.91dcba358223f4ada572a60e999c28132620d6ae.de-de.xlf"
$de.Range("K2").Value = "2016-10-24 09:53:17"

# Row 3 : de4b381b-7123-43c1-a4ab-173c364b5d43
$de.Range("I3").Value = "de4b381b-7123-43c1-a4ab-173c364b5d43.md"
$de.Hyperlinks.Add($de.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/64e54dce710ff6806c33fcbd5c6c5c24cea91328/e2e/de4b381b-7123-43c1-a4ab-173c364b5d43.md", "", "", "de4b381b-7123-43c1-a4ab-173c364b5d43.md") | Out-Null
$de.Range("J3").Value = "de4b381b-7123-43c1-a4ab-173c364b5d43.9c06a25dcb376ef7f97389c27a3adf04fa2ec2ca.de-de.xlf"
$de.Range("K3").Value = "2016-10-24 09:53:17"

$de.Columns.Item(3).ColumnWidth = 29.14   # raw width 30.0 (closest achievable to target 29.9777050018311)
$de.Columns.Item(9).ColumnWidth = 39.14   # raw width 40.0
$de.Columns.Item(10).ColumnWidth = 39.14  # raw width 40.0
